# Applies the row-rotation edit described by the diff:
#   row 45 <- (old) row 47 data, row 46 <- (old) row 45 data, row 47 <- (old) row 46 data
#   row 55 <- (old) row 56 data, row 56 <- (old) row 57 data, row 57 <- (old) row 55 data
# Each destination cell is only written when its value actually changes, so
# untouched cells keep their original formatting/encoding.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45 gets the data previously in row 47 ---
$ws.Range("A45").Value = 111529093
$ws.Range("B45").Value = 82949
$ws.Range("D45").Value = "NT"
$ws.Range("E45").Value = 5589
$ws.Range("F45").Value = "Rödbrun klubbdyna"
$ws.Range("G45").Value = "Trichoderma nybergianum"
$ws.Range("H45").Value = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
$ws.Range("P45").Value = "Bye kalkbarrskogs naturreservat (Bye kalkbarrskogs naturreservat), Jmt"
$ws.Range("Q45").Value = 485442.3376153786
$ws.Range("R45").Value = 6995847.791586295
$ws.Range("Y45").NumberFormat = "@"
$ws.Range("Y45").Value = "2023-08-17"
$ws.Range("Z45").Value = "00:00"
$ws.Range("AA45").NumberFormat = "@"
$ws.Range("AA45").Value = "2023-08-17"
$ws.Range("AB45").Value = "00:00"
$ws.Range("AW45").Value = "Rashid Kadhim"
$ws.Range("AX45").Value = "Rashid Kadhim, Thomas Stålhandske, Hugo Ström"

# --- Row 46 gets the data previously in row 45 ---
$ws.Range("A46").Value = 111646113
$ws.Range("B46").Value = 87995
$ws.Range("E46").Value = 1594
$ws.Range("F46").Value = "Svartfjällig musseron"
$ws.Range("G46").Value = "Tricholoma atrosquamosum"
$ws.Range("H46").Value = "Sacc."
$ws.Range("Q46").Value = 485314.111892351
$ws.Range("R46").Value = 6995879.171382442
$ws.Range("Z46").Value = "13:42"
$ws.Range("AB46").Value = "13:42"

# --- Row 47 gets the data previously in row 46 ---
$ws.Range("A47").Value = 111645939
$ws.Range("B47").Value = 88956
$ws.Range("D47").Value = "VU"
$ws.Range("E47").Value = 5747
$ws.Range("F47").Value = "Läderdoftande fingersvamp"
$ws.Range("G47").Value = "Ramaria safraniolens"
$ws.Range("H47").Value = "Christian"
$ws.Range("P47").Value = "Källmyren (Källmyren), Jmt"
$ws.Range("Q47").Value = 485250.2046207946
$ws.Range("R47").Value = 6995790.225604231
$ws.Range("Y47").NumberFormat = "@"
$ws.Range("Y47").Value = "2023-08-23"
$ws.Range("Z47").Value = "13:30"
$ws.Range("AA47").NumberFormat = "@"
$ws.Range("AA47").Value = "2023-08-23"
$ws.Range("AB47").Value = "13:30"
$ws.Range("AW47").Value = "Andreas Öster"
$ws.Range("AX47").Value = "Andreas Öster"

# --- Row 55 gets the data previously in row 56 ---
$ws.Range("A55").Value = 111646292
$ws.Range("B55").Value = 88946
$ws.Range("E55").Value = 256335
$ws.Range("F55").Value = "Taggfingersvamp"
$ws.Range("G55").Value = "Ramaria karstenii"
$ws.Range("H55").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("P55").Value = "Bye (Bye), Jmt"
$ws.Range("Q55").Value = 485380.6216548251
$ws.Range("R55").Value = 6995741.1661748
$ws.Range("Z55").Value = "00:00"
$ws.Range("AB55").Value = "00:00"
$ws.Range("AW55").Value = "Rashid Kadhim"
$ws.Range("AX55").Value = "Rashid Kadhim"

# --- Row 56 gets the data previously in row 57 ---
$ws.Range("A56").Value = 111645934
$ws.Range("B56").Value = 88915
$ws.Range("D56").Value = "NT"
$ws.Range("E56").Value = 5734
$ws.Range("F56").Value = "Druvfingersvamp"
$ws.Range("G56").Value = "Ramaria botrytis"
$ws.Range("H56").Value = "(Pers.:Fr.) Bourdot"
$ws.Range("P56").Value = "Källmyren (Källmyren), Jmt"
$ws.Range("Q56").Value = 485252.940604815
$ws.Range("R56").Value = 6995793.384020397
$ws.Range("Z56").Value = "13:30"
$ws.Range("AB56").Value = "13:30"
$ws.Range("AW56").Value = "Andreas Öster"
$ws.Range("AX56").Value = "Andreas Öster"

# --- Row 57 gets the data previously in row 55 ---
$ws.Range("A57").Value = 111645786
$ws.Range("B57").Value = 88956
$ws.Range("D57").Value = "VU"
$ws.Range("E57").Value = 5747
$ws.Range("F57").Value = "Läderdoftande fingersvamp"
$ws.Range("G57").Value = "Ramaria safraniolens"
$ws.Range("H57").Value = "Christian"
$ws.Range("Q57").Value = 485330.9609580904
$ws.Range("R57").Value = 6995793.48329893
$ws.Range("Z57").Value = "13:21"
$ws.Range("AB57").Value = "13:21"

